$d = $word.ActiveDocument

# 1) "Éttermek" mezőlista bővítése: férőhely,leírás,státusz -> férőhely,házhozszállítás, leírás,tipus,wifi,státusz,kép,weboldal,facebook,nyitvavane
$d.Content.Find.Execute(
    "férőhely,leírás,státusz)", $false, $false, $false, $false, $false,
    $true, 1, $false,
    "férőhely,házhozszállítás, leírás,tipus,wifi,státusz,kép,weboldal,facebook,nyitvavane)",
    2)

# 2) "Nyitvatartás" mezőlista bővítése: nap, nyitás -> nap,napid, nyitás
$d.Content.Find.Execute(
    ", nap, nyitás, zárás)", $false, $false, $false, $false, $false,
    $true, 1, $false,
    ", nap,napid, nyitás, zárás)",
    2)

# 3) "Helyfoglalás" mezőlista egyszerűsítése: kezdés ideje,végzés ideje,fő -> kezdés ideje,fő
$d.Content.Find.Execute(
    "kezdés ideje,végzés ideje,fő)", $false, $false, $false, $false, $false,
    $true, 1, $false,
    "kezdés ideje,fő)",
    2)

# 4) Felesleges "_GoBack" bookmark eltávolítása
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
